$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8855.637000000001
$ws.Range("I86").Value = 6666.6665
$ws.Range("J86").Value = 9676.5
$ws.Range("K86").Value = 6666.6665
$ws.Range("L86").Value = 9676.5
$ws.Range("M86").Value = -5543.6665
$ws.Range("N86").Value = -11922.5
$ws.Range("H89").Value = 8855.637000000001
$ws.Range("I89").Value = 6666.6665
$ws.Range("J89").Value = 9676.5
$ws.Range("K89").Value = 33333.3325
$ws.Range("L89").Value = 48382.5
$ws.Range("M89").Value = -27717.3325
$ws.Range("N89").Value = -59614.5
$ws.Range("H92").Value = 1556.2858
$ws.Range("I92").Value = 1399
$ws.Range("J92").Value = 2500
$ws.Range("K92").Value = 1399
$ws.Range("L92").Value = 2500
$ws.Range("M92").Value = -151
$ws.Range("N92").Value = -4996
$ws.Range("H100").Value = 1996.5333
$ws.Range("I100").Value = 1709.6
$ws.Range("J100").Value = 2140
$ws.Range("K100").Value = 1709.6
$ws.Range("L100").Value = 2140
$ws.Range("M100").Value = -1168.6
$ws.Range("N100").Value = -3222
$ws.Range("H132").Value = 2038.762
$ws.Range("I132").Value = 2359.9333
$ws.Range("J132").Value = 1235.8334
$ws.Range("K132").Value = 7079.7999
$ws.Range("L132").Value = 3707.5002
$ws.Range("M132").Value = -4549.7999
$ws.Range("N132").Value = -8767.5002
$ws.Range("H137").Value = 1917637.9
$ws.Range("I137").Value = 1733.1
$ws.Range("J137").Value = 6175204
$ws.Range("K137").Value = 5199.299999999999
$ws.Range("L137").Value = 18525612
$ws.Range("M137").Value = -2649.299999999999
$ws.Range("N137").Value = -18530712
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1163.8
$ws.Range("I97").Value = 1119.7142
$ws.Range("J97").Value = 1266.6666
$ws.Range("K97").Value = 1119.7142
$ws.Range("L97").Value = 1266.6666
$ws.Range("M97").Value = -623.7141999999999
$ws.Range("N97").Value = -2258.6666
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 769
$ws.Range("I94").Value = 690.6667
$ws.Range("J94").Value = 910
$ws.Range("K94").Value = 690.6667
$ws.Range("L94").Value = 910
$ws.Range("M94").Value = -239.6667
$ws.Range("N94").Value = -1812
$ws.Range("H99").Value = 2213.6667
$ws.Range("I99").Value = 1806.7142
$ws.Range("J99").Value = 2472.6365
$ws.Range("K99").Value = 1806.7142
$ws.Range("L99").Value = 2472.6365
$ws.Range("M99").Value = -308.7141999999999
$ws.Range("N99").Value = -5468.636500000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1709.6471
$ws.Range("I58").Value = 1670.9333
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1670.9333
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -1467.9333
$ws.Range("N58").Value = -2406
$ws.Range("H99").Value = 4209.6665
$ws.Range("I99").Value = 4326
$ws.Range("J99").Value = 4093.3333
$ws.Range("K99").Value = 4326
$ws.Range("L99").Value = 4093.3333
$ws.Range("M99").Value = -2828
$ws.Range("N99").Value = -7089.3333
$ws.Range("H126").Value = 4209.6665
$ws.Range("I126").Value = 4326
$ws.Range("J126").Value = 4093.3333
$ws.Range("K126").Value = 12978
$ws.Range("L126").Value = 12279.9999
$ws.Range("M126").Value = -10508
$ws.Range("N126").Value = -17219.9999
$ws.Range("H132").Value = 2471.8
$ws.Range("I132").Value = 2177.375
$ws.Range("J132").Value = 3649.5
$ws.Range("K132").Value = 6532.125
$ws.Range("L132").Value = 10948.5
$ws.Range("M132").Value = -4002.125
$ws.Range("N132").Value = -16008.5
$ws.Range("H134").Value = 1893.1904
$ws.Range("I134").Value = 1816.0588
$ws.Range("K134").Value = 5448.1764
$ws.Range("M134").Value = -2913.1764
$ws.Range("H136").Value = 1709.6471
$ws.Range("I136").Value = 1670.9333
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5012.7999
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2462.7999
$ws.Range("N136").Value = -11100
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 82.85714
$ws.Range("I8").Value = 82.85714
$ws.Range("K8").Value = 248.57142
$ws.Range("M8").Value = -109.57142
$ws.Range("H80").Value = 1634
$ws.Range("J80").Value = 1608.75
$ws.Range("L80").Value = 4826.25
$ws.Range("N80").Value = -6698.25
$ws.Range("H83").Value = 1634
$ws.Range("J83").Value = 1608.75
$ws.Range("L83").Value = 14478.75
$ws.Range("N83").Value = -23838.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2755.6
$ws.Range("I80").Value = 2659.0476
$ws.Range("J80").Value = 3262.5
$ws.Range("K80").Value = 2659.0476
$ws.Range("L80").Value = 3262.5
$ws.Range("M80").Value = -1661.0476
$ws.Range("N80").Value = -5258.5
$ws.Range("H83").Value = 2755.6
$ws.Range("I83").Value = 2659.0476
$ws.Range("J83").Value = 3262.5
$ws.Range("K83").Value = 13295.238
$ws.Range("L83").Value = 16312.5
$ws.Range("M83").Value = -8303.237999999999
$ws.Range("N83").Value = -26296.5
$ws.Range("H107").Value = 1435.4286
$ws.Range("I107").Value = 1416
$ws.Range("J107").Value = 1450
$ws.Range("K107").Value = 1416
$ws.Range("L107").Value = 1450
$ws.Range("M107").Value = 504
$ws.Range("N107").Value = -5290
$ws.Range("H126").Value = 2156.074
$ws.Range("I126").Value = 1327.5714
$ws.Range("J126").Value = 3048.3076
$ws.Range("K126").Value = 3982.7142
$ws.Range("L126").Value = 9144.9228
$ws.Range("M126").Value = -1512.7142
$ws.Range("N126").Value = -14084.9228
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2463.4211
$ws.Range("I82").Value = 3002
$ws.Range("J82").Value = 2433.5
$ws.Range("K82").Value = 3002
$ws.Range("L82").Value = 2433.5
$ws.Range("M82").Value = -2641
$ws.Range("N82").Value = -3155.5
$ws.Range("H85").Value = 2463.4211
$ws.Range("I85").Value = 3002
$ws.Range("J85").Value = 2433.5
$ws.Range("K85").Value = 3002
$ws.Range("L85").Value = 2433.5
$ws.Range("M85").Value = -1754
$ws.Range("N85").Value = -4929.5
$ws.Range("H93").Value = 33115.547
$ws.Range("I93").Value = 1495.9375
$ws.Range("J93").Value = 117434.5
$ws.Range("K93").Value = 1495.9375
$ws.Range("L93").Value = 117434.5
$ws.Range("M93").Value = -247.9375
$ws.Range("N93").Value = -119930.5
$ws.Range("H132").Value = 4131
$ws.Range("I132").Value = 4044.9443
$ws.Range("J132").Value = 4324.625
$ws.Range("K132").Value = 12134.8329
$ws.Range("L132").Value = 12973.875
$ws.Range("M132").Value = -9604.832900000001
$ws.Range("N132").Value = -18033.875
$ws.Range("H136").Value = 2449.25
$ws.Range("I136").Value = 2535.0908
$ws.Range("J136").Value = 1505
$ws.Range("K136").Value = 7605.2724
$ws.Range("L136").Value = 4515
$ws.Range("M136").Value = -5055.2724
$ws.Range("N136").Value = -9615
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1523.5333
$ws.Range("I96").Value = 1178.4166
$ws.Range("J96").Value = 2904
$ws.Range("K96").Value = 1178.4166
$ws.Range("L96").Value = 2904
$ws.Range("M96").Value = 194.5834
$ws.Range("N96").Value = -5650
$ws.Range("H100").Value = 1120.5883
$ws.Range("I100").Value = 1120.5883
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2241.1766
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1700.1766
$ws.Range("N100").ClearContents()
$ws.Range("H126").Value = 1251.9722
$ws.Range("I126").Value = 716.8095
$ws.Range("J126").Value = 2001.2
$ws.Range("K126").Value = 2150.4285
$ws.Range("L126").Value = 6003.6
$ws.Range("M126").Value = 319.5715
$ws.Range("N126").Value = -10943.6
$ws.Range("H132").Value = 3249.842
$ws.Range("I132").Value = 3242.3572
$ws.Range("J132").Value = 3270.8
$ws.Range("K132").Value = 9727.071599999999
$ws.Range("L132").Value = 9812.400000000001
$ws.Range("M132").Value = -7197.071599999999
$ws.Range("N132").Value = -14872.4
